$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.474.59"
Set-TextValue "E2" "  +1.50%  "

Set-TextValue "D3" "1.678.76"
Set-TextValue "E3" "  +2.76%  "

Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  +0.03%  "

Set-TextValue "D5" "216.71"
Set-TextValue "E5" "  +1.23%  "

Set-TextValue "D6" "0.5303"
Set-TextValue "E6" "  +1.66%  "

Set-TextValue "D7" "1.002"
Set-TextValue "E7" "  +0.02%  "

Set-TextValue "D8" "0.2699"
Set-TextValue "E8" "  +3.79%  "

Set-TextValue "D9" "0.06402"
Set-TextValue "E9" "  +1.84%  "

Set-TextValue "D10" "21.72"
Set-TextValue "E10" "  +5.58%  "

Set-TextValue "D11" "0.07813"
Set-TextValue "E11" "  +2.61%  "

Set-TextValue "D12" "1.684.66"
Set-TextValue "E12" "  +3.16%  "

Set-TextValue "D13" "4.501"
Set-TextValue "E13" "  +1.95%  "

Set-TextValue "D14" "0.5563"
Set-TextValue "E14" "  +1.29%  "

Set-TextValue "D15" "0.0₅8324"
Set-TextValue "E15" "  +3.26%  "

Set-TextValue "D16" "65.57"
Set-TextValue "E16" "  +1.22%  "

Set-TextValue "D17" "26.526.43"
Set-TextValue "E17" "  +1.82%  "

Set-TextValue "D18" "1.002"
Set-TextValue "E18" "  -0.02%  "

Set-TextValue "D19" "4.741"
Set-TextValue "E19" "  +1.50%  "

Set-TextValue "D20" "193.83"
Set-TextValue "E20" "  +2.97%  "

Set-TextValue "D21" "10.34"
Set-TextValue "E21" "  +1.93%  "

Set-TextValue "D22" "6.338"
Set-TextValue "E22" "  +3.39%  "

Set-TextValue "D23" "1.003"
Set-TextValue "E23" "  +0.03%  "

Set-TextValue "D24" "142.53"
Set-TextValue "E24" "  -2.13%  "

Set-TextValue "D25" "0.1286"
Set-TextValue "E25" "  +5.83%  "

Set-TextValue "D26" "7.416"
Set-TextValue "E26" "  +0.27%  "

Set-TextValue "D27" "16.27"
Set-TextValue "E27" "  +2.74%  "

Set-TextValue "D28" "1.421"
Set-TextValue "E28" "  +2.16%  "

Set-TextValue "D29" "0.06227"
Set-TextValue "E29" "  +6.28%  "

Set-TextValue "D30" "1.270"

Set-TextValue "D31" "3.608"
Set-TextValue "E31" "  +5.21%  "

Set-TextValue "D32" "3.445"
Set-TextValue "E32" "  +1.37%  "

Set-TextValue "D33" "1.676"
Set-TextValue "E33" "  +2.40%  "

Set-TextValue "D34" "1.008"
Set-TextValue "E34" "  +2.54%  "

Set-TextValue "E35" "  +1.39%  "

Set-TextValue "D36" "2.783"
Set-TextValue "E36" "  +0.65%  "

Set-TextValue "D37" "0.6025"
Set-TextValue "E37" "  +5.00%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D38" "6.178"
Set-TextValue "E38" "  +9.23%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.01634"
Set-TextValue "E39" "  +1.20%  "

Set-TextValue "D40" "1.083.97"
Set-TextValue "E40" "  +4.94%  "

Set-TextValue "D41" "0.8626"
Set-TextValue "E41" "  +0.70%  "

Set-TextValue "D42" "1.001"
Set-TextValue "E42" "  -0.02%  "

Set-TextValue "D43" "100.21"
Set-TextValue "E43" "  -0.06%  "

Set-TextValue "D44" "1.824.78"
Set-TextValue "E44" "  +2.32%  "

Set-TextValue "D45" "0.0₈110"
Set-TextValue "E45" "  +3.23%  "

Set-TextValue "D46" "57.14"
Set-TextValue "E46" "  +3.27%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D47" "1.003"
Set-TextValue "E47" "  +0.64%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "8.125"
Set-TextValue "E48" "  +0.93%  "

Set-TextValue "E49" "  +0.80%  "

Set-TextValue "D50" "6.034"
Set-TextValue "E50" "  +2.82%  "

Set-TextValue "E51" "  +0.38%  "
